$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# xlEdgeBottom = 9 ; xlLineStyleNone = -4142
$xlEdgeBottom = 9
$xlLineStyleNone = -4142

# Row 18 (Sectioning): mark as "fully implemented" instead of "not implemented".
# This row is no longer the last one of its bordered block, so drop the
# bottom border (matches the other "interior" rows of the block, e.g. B13-B16).
$ws.Range("B18").Value = "fully implemented"
$ws.Range("B18").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone

# Row 20 (Styling): mark as "fully implemented" instead of "not implemented"
$ws.Range("B20").Value = "fully implemented"
$ws.Range("B20").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone

# Row 22 (Cross-Validation): mark as "partially implemented/limited support)" and add a comment
$ws.Range("B22").Value = "partially implemented/limited support)"
$ws.Range("B22").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone
$ws.Range("C22").Value = "links to QL are enforced, however some more validation should be implemented"

# Update the current selection, matching the saved cursor position
$ws.Range("C23").Select()
